$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5505.467
$ws.Range("I62").Value = 4512.143
$ws.Range("J62").Value = 6374.625
$ws.Range("K62").Value = 4512.143
$ws.Range("L62").Value = 6374.625
$ws.Range("M62").Value = -3888.143
$ws.Range("N62").Value = -7622.625
$ws.Range("H65").Value = 5505.467
$ws.Range("I65").Value = 4512.143
$ws.Range("J65").Value = 6374.625
$ws.Range("K65").Value = 22560.715
$ws.Range("L65").Value = 31873.125
$ws.Range("M65").Value = -19440.715
$ws.Range("N65").Value = -38113.125
$ws.Range("H86").Value = 10239.091
$ws.Range("I86").Value = 1211
$ws.Range("K86").Value = 1211
$ws.Range("M86").Value = -88
$ws.Range("H89").Value = 10239.091
$ws.Range("I89").Value = 1211
$ws.Range("K89").Value = 6055
$ws.Range("M89").Value = -439
$ws.Range("H106").Value = 6291215.5
$ws.Range("I106").Value = 8773149
$ws.Range("J106").Value = 3651.7334
$ws.Range("K106").Value = 8773149
$ws.Range("L106").Value = 3651.7334
$ws.Range("M106").Value = -8772518
$ws.Range("N106").Value = -4913.7334
$ws.Range("H112").Value = 2850132
$ws.Range("J112").Value = 2850132
$ws.Range("L112").Value = 8550396
$ws.Range("N112").Value = -8552612
$ws.Range("H113").Value = 71434104
$ws.Range("J113").Value = 7800
$ws.Range("L113").Value = 7800
$ws.Range("N113").Value = -14308
$ws.Range("H116").Value = 4366.952
$ws.Range("I116").Value = 2128.5715
$ws.Range("J116").Value = 5486.143
$ws.Range("K116").Value = 2128.5715
$ws.Range("L116").Value = 5486.143
$ws.Range("M116").Value = 1313.4285
$ws.Range("N116").Value = -12370.143
$ws.Range("H137").Value = 1884.6923
$ws.Range("I137").Value = 1649.8
$ws.Range("K137").Value = 4949.4
$ws.Range("M137").Value = -2399.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3560.6667
$ws.Range("I32").Value = 2807.908
$ws.Range("K32").Value = 2807.908
$ws.Range("M32").Value = -2520.908
$ws.Range("H45").Value = 2569.8667
$ws.Range("I45").Value = 2351
$ws.Range("J45").Value = 2820
$ws.Range("K45").Value = 2351
$ws.Range("L45").Value = 2820
$ws.Range("M45").Value = -1974
$ws.Range("N45").Value = -3574
$ws.Range("H61").Value = 2135.4285
$ws.Range("I61").Value = 1595.5
$ws.Range("K61").Value = 1595.5
$ws.Range("M61").Value = -1383.5
$ws.Range("H74").Value = 50001896
$ws.Range("I74").Value = 58824584
$ws.Range("K74").Value = 58824584
$ws.Range("M74").Value = -58823710
$ws.Range("H77").Value = 50001896
$ws.Range("I77").Value = 58824584
$ws.Range("K77").Value = 294122920
$ws.Range("M77").Value = -294118552
$ws.Range("H110").Value = 1178.5
$ws.Range("I110").Value = 1105
$ws.Range("J110").Value = 1399
$ws.Range("K110").Value = 1105
$ws.Range("L110").Value = 1399
$ws.Range("M110").Value = 940
$ws.Range("N110").Value = -5489
$ws.Range("H122").Value = 2306.5144
$ws.Range("I122").Value = 1954.875
$ws.Range("J122").Value = 3073.7273
$ws.Range("K122").Value = 5864.625
$ws.Range("L122").Value = 9221.1819
$ws.Range("M122").Value = -3414.625
$ws.Range("N122").Value = -14121.1819
$ws.Range("H132").Value = 16968.455
$ws.Range("I132").Value = 1464
$ws.Range("K132").Value = 4392
$ws.Range("M132").Value = -1862
$ws.Range("H136").Value = 2135.4285
$ws.Range("I136").Value = 1595.5
$ws.Range("K136").Value = 4786.5
$ws.Range("M136").Value = -2236.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1551.3414
$ws.Range("I86").Value = 1432.2
$ws.Range("J86").Value = 1737.5
$ws.Range("K86").Value = 1432.2
$ws.Range("L86").Value = 1737.5
$ws.Range("M86").Value = -309.2
$ws.Range("N86").Value = -3983.5
$ws.Range("H89").Value = 1551.3414
$ws.Range("I89").Value = 1432.2
$ws.Range("J89").Value = 1737.5
$ws.Range("K89").Value = 7161
$ws.Range("L89").Value = 8687.5
$ws.Range("M89").Value = -1545
$ws.Range("N89").Value = -19919.5
$ws.Range("H107").Value = 653.3333
$ws.Range("I107").Value = 684
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 684
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1236
$ws.Range("N107").Value = -4340

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 748.3333
$ws.Range("I12").Value = 230
$ws.Range("K12").Value = 230
$ws.Range("M12").Value = -60
$ws.Range("H58").Value = 33338.875
$ws.Range("I58").Value = 2798.8
$ws.Range("J58").Value = 47220.727
$ws.Range("K58").Value = 2798.8
$ws.Range("L58").Value = 47220.727
$ws.Range("M58").Value = -2595.8
$ws.Range("N58").Value = -47626.727
$ws.Range("H105").Value = 31251130
$ws.Range("I105").Value = 125000000
$ws.Range("J105").Value = 1507.3334
$ws.Range("K105").Value = 125000000
$ws.Range("L105").Value = 1507.3334
$ws.Range("M105").Value = -124998253
$ws.Range("N105").Value = -5001.3334
$ws.Range("H132").Value = 2877.3333
$ws.Range("I132").Value = 2274.8
$ws.Range("J132").Value = 4598.857
$ws.Range("K132").Value = 6824.400000000001
$ws.Range("L132").Value = 13796.571
$ws.Range("M132").Value = -4294.400000000001
$ws.Range("N132").Value = -18856.571
$ws.Range("H136").Value = 33338.875
$ws.Range("I136").Value = 2798.8
$ws.Range("J136").Value = 47220.727
$ws.Range("K136").Value = 8396.400000000001
$ws.Range("L136").Value = 141662.181
$ws.Range("M136").Value = -5846.400000000001
$ws.Range("N136").Value = -146762.181

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 735.37
$ws.Range("I131").Value = 357.66666
$ws.Range("J131").Value = 772.7252999999999
$ws.Range("K131").Value = 1072.99998
$ws.Range("L131").Value = 2318.1759
$ws.Range("M131").Value = 3967.00002
$ws.Range("N131").Value = -12398.1759

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15937.5
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 17928.572
$ws.Range("K70").Value = 2000
$ws.Range("L70").Value = 17928.572
$ws.Range("M70").Value = -1730
$ws.Range("N70").Value = -18468.572
$ws.Range("H73").Value = 15937.5
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 17928.572
$ws.Range("K73").Value = 2000
$ws.Range("L73").Value = 17928.572
$ws.Range("M73").Value = -1064
$ws.Range("N73").Value = -19800.572
$ws.Range("H102").Value = 16668992
$ws.Range("J102").Value = 2782.8
$ws.Range("L102").Value = 2782.8
$ws.Range("N102").Value = -6026.8
$ws.Range("H122").Value = 88890700
$ws.Range("I122").Value = 30304866
$ws.Range("K122").Value = 90914598
$ws.Range("M122").Value = -90912148

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2879.9524
$ws.Range("I40").Value = 2088.9
$ws.Range("J40").Value = 3599.0908
$ws.Range("K40").Value = 2088.9
$ws.Range("L40").Value = 3599.0908
$ws.Range("M40").Value = -1952.9
$ws.Range("N40").Value = -3871.0908
$ws.Range("H122").Value = 983230.1
$ws.Range("I122").Value = 1510784.9
$ws.Range("J122").Value = 3485.5715
$ws.Range("K122").Value = 4532354.699999999
$ws.Range("L122").Value = 10456.7145
$ws.Range("M122").Value = -4529904.699999999
$ws.Range("N122").Value = -15356.7145
$ws.Range("H132").Value = 929572.5600000001
$ws.Range("I132").Value = 2412109
$ws.Range("J132").Value = 2987.25
$ws.Range("K132").Value = 7236327
$ws.Range("L132").Value = 8961.75
$ws.Range("M132").Value = -7233797
$ws.Range("N132").Value = -14021.75
$ws.Range("H136").Value = 1888.8889
$ws.Range("I136").Value = 1888.8889
$ws.Range("K136").Value = 5666.6667
$ws.Range("M136").Value = -3116.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1733.3334
$ws.Range("I132").Value = 1425
$ws.Range("J132").Value = 1980
$ws.Range("K132").Value = 4275
$ws.Range("L132").Value = 5940
$ws.Range("M132").Value = -1745
$ws.Range("N132").Value = -11000
$ws.Range("H136").Value = 31282624
$ws.Range("I136").Value = 46922476
$ws.Range("J136").Value = 2919.0908
$ws.Range("K136").Value = 140767428
$ws.Range("L136").Value = 8757.2724
$ws.Range("M136").Value = -140764878
$ws.Range("N136").Value = -13857.2724
